$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the special underline-style formatting from rows 7-14 (back to default "Normal" style)
$ws.Range("A7:F14").ClearFormats()

# Remove the stray styled-only cell I7 and clear the leftover styled D16 cell
$ws.Range("I7").Clear()
$ws.Range("D16").Clear()

# New employee rows 15-21
$ws.Range("A15").Value = 5124
$ws.Range("B15").Value = "Igor"
$ws.Range("C15").Value = "Banaa"
$ws.Range("D15").Value = 5532525325
$ws.Range("E15").Value = "Teste"
$ws.Range("F15").Value = "Malaquias"

$ws.Range("A16").Value = 45543
$ws.Range("B16").Value = "Marcos"
$ws.Range("C16").Value = "Polo"
$ws.Range("D16").Value = 41424214214
$ws.Range("E16").Value = "Agua"
$ws.Range("F16").Value = "Arroz"

$ws.Range("A17").Value = 7895
$ws.Range("B17").Value = "Jeferson"
$ws.Range("C17").Value = "Irineu"
$ws.Range("D17").Value = 422526236
$ws.Range("E17").Value = "Banana"
$ws.Range("F17").Value = "Nabo"

$ws.Range("A18").Value = 46755
$ws.Range("B18").Value = "Gabriel"
$ws.Range("C18").Value = "Jesus"
$ws.Range("D18").Value = 526752532
$ws.Range("E18").Value = "Banana"
$ws.Range("F18").Value = "Nabo"

$ws.Range("A19").Value = 31465
$ws.Range("B19").Value = "Igor"
$ws.Range("C19").Value = "Banaa"
$ws.Range("D19").Value = 21467534732
$ws.Range("E19").Value = "Teste"
$ws.Range("F19").Value = "Malaquias"

$ws.Range("A20").Value = 6788
$ws.Range("B20").Value = "Marcos"
$ws.Range("C20").Value = "Polo"
$ws.Range("D20").Value = 4414124214
$ws.Range("E20").Value = "Agua"
$ws.Range("F20").Value = "Arroz"

$ws.Range("A21").Value = 21342
$ws.Range("B21").Value = "Renato"
$ws.Range("C21").Value = "Kenzo"
$ws.Range("D21").Value = 5511999381877
$ws.Range("E21").Value = "Teste"
$ws.Range("F21").Value = "Banana"

$ws.Range("I17").Select()
